# Applies the commit "Update CDA Logical model for ST.r2b" to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 2 rename: "Include from RoleClass" -> "Include #0" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- Sheet 1 (Metadata) updates ---
$ws1 = $wb.Worksheets.Item(1)

# Version value update
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Date value update
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for the new "Jurisdiction" property,
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws1.Rows.Item(11).Insert()

# Copy formatting from the row above (Contact, row 10) onto the newly inserted row
# so the new row keeps the same cell style (s="2") as the rest of the data rows.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" row
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# The "Description" row (now row 12 after the insert) gets its Value populated
$ws1.Range("A12").Value = "Description"
$ws1.Range("B12").Value = "Used to represent the role(s) of those who should receive a copy of a document - limited to values allowed in original CDA definition"

Write-Output "edit applied"
